# Update two-digit multiplication problem/answer strings in the table.
$d = $word.ActiveDocument

$pairs = @(
    @{old = "22×82=1804"; new = "15×38=570"},
    @{old = "27×87=2349"; new = "75×37=2775"},
    @{old = "19×42=798";  new = "66×32=2112"},
    @{old = "88×32=2816"; new = "21×95=1995"},
    @{old = "91×65=5915"; new = "14×61=854"},
    @{old = "28×68=1904"; new = "41×42=1722"},
    @{old = "47×51=2397"; new = "62×29=1798"},
    @{old = "70×27=1890"; new = "62×60=3720"},
    @{old = "82×71=5822"; new = "83×84=6972"},
    @{old = "49×86=4214"; new = "74×60=4440"},
    @{old = "79×82=6478"; new = "57×19=1083"},
    @{old = "70×87=6090"; new = "61×52=3172"},
    @{old = "74×82=6068"; new = "69×20=1380"},
    @{old = "45×58=2610"; new = "12×98=1176"},
    @{old = "74×88=6512"; new = "55×45=2475"},
    @{old = "55×48=2640"; new = "16×44=704"},
    @{old = "59×78=4602"; new = "51×15=765"},
    @{old = "58×56=3248"; new = "41×60=2460"},
    @{old = "58×58=3364"; new = "86×87=7482"},
    @{old = "77×32=2464"; new = "80×29=2320"},
    @{old = "72×36=2592"; new = "83×28=2324"},
    @{old = "33×89=2937"; new = "56×48=2688"},
    @{old = "24×94=2256"; new = "70×99=6930"},
    @{old = "53×97=5141"; new = "83×23=1909"},
    @{old = "41×69=2829"; new = "71×22=1562"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
